# Fruta / hortaliza, semanal
# Insert 3 new weekly records at the top of the data block (rows 782-784),
# pushing the existing data down by 3 rows (dimension grows from T857 to T860).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows right before the current row 782.
$ws.Rows("782:784").Insert()

# Row 782: Packham's Triumph - Especial
$ws.Cells.Item(782, 1).Value = 9
$ws.Cells.Item(782, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(782, 3).Value = "Metropolitana"
$ws.Cells.Item(782, 4).Value = 44769
$ws.Cells.Item(782, 5).Value = 13
$ws.Cells.Item(782, 6).Value = "Fruta"
$ws.Cells.Item(782, 7).Value = 100104
$ws.Cells.Item(782, 8).Value = "Frutos de pepita"
$ws.Cells.Item(782, 9).Value = 100104005
$ws.Cells.Item(782, 10).Value = "Pera"
$ws.Cells.Item(782, 11).Value = "Packham's Triumph"
$ws.Cells.Item(782, 12).Value = "Especial"
$ws.Cells.Item(782, 13).Value = 4
$ws.Cells.Item(782, 14).Value = 160000
$ws.Cells.Item(782, 15).Value = 160000
$ws.Cells.Item(782, 16).Value = 160000
$ws.Cells.Item(782, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(782, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(782, 19).Value = 356
$ws.Cells.Item(782, 20).Value = 450

# Row 783: Packham's Triumph - Primera
$ws.Cells.Item(783, 1).Value = 9
$ws.Cells.Item(783, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(783, 3).Value = "Metropolitana"
$ws.Cells.Item(783, 4).Value = 44769
$ws.Cells.Item(783, 5).Value = 13
$ws.Cells.Item(783, 6).Value = "Fruta"
$ws.Cells.Item(783, 7).Value = 100104
$ws.Cells.Item(783, 8).Value = "Frutos de pepita"
$ws.Cells.Item(783, 9).Value = 100104005
$ws.Cells.Item(783, 10).Value = "Pera"
$ws.Cells.Item(783, 11).Value = "Packham's Triumph"
$ws.Cells.Item(783, 12).Value = "Primera"
$ws.Cells.Item(783, 13).Value = 8
$ws.Cells.Item(783, 14).Value = 150000
$ws.Cells.Item(783, 15).Value = 150000
$ws.Cells.Item(783, 16).Value = 150000
$ws.Cells.Item(783, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(783, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(783, 19).Value = 333
$ws.Cells.Item(783, 20).Value = 450

# Row 784: Packham's Triumph - Segunda
$ws.Cells.Item(784, 1).Value = 9
$ws.Cells.Item(784, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(784, 3).Value = "Metropolitana"
$ws.Cells.Item(784, 4).Value = 44769
$ws.Cells.Item(784, 5).Value = 13
$ws.Cells.Item(784, 6).Value = "Fruta"
$ws.Cells.Item(784, 7).Value = 100104
$ws.Cells.Item(784, 8).Value = "Frutos de pepita"
$ws.Cells.Item(784, 9).Value = 100104005
$ws.Cells.Item(784, 10).Value = "Pera"
$ws.Cells.Item(784, 11).Value = "Packham's Triumph"
$ws.Cells.Item(784, 12).Value = "Segunda"
$ws.Cells.Item(784, 13).Value = 10
$ws.Cells.Item(784, 14).Value = 130000
$ws.Cells.Item(784, 15).Value = 130000
$ws.Cells.Item(784, 16).Value = 130000
$ws.Cells.Item(784, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(784, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(784, 19).Value = 289
$ws.Cells.Item(784, 20).Value = 450
